# Apply updated crypto price (D) and 1h volume change (E) values.
# Column D values that look like plain numbers are written with a leading
# apostrophe (text qualifier) so Excel keeps/stores them as text -- exactly
# like the source data -- instead of coercing to a Number and silently
# dropping significant trailing zeros (e.g. "8.90" -> 8.9).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.792.77"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "3.365.85"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'568.74"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'137.74"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "'7.65"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "'0.381"
$ws.Range("E11").Value = "  -4.27%  "
$ws.Range("D12").Value = "3.939.44"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "'27.72"
$ws.Range("E14").Value = "  -2.11%  "
$ws.Range("D15").Value = "3.349.78"
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "60.906.96"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  -2.58%  "
$ws.Range("D19").Value = "'13.51"
$ws.Range("E19").Value = "  -3.73%  "
$ws.Range("D20").Value = "'8.90"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "'381.67"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").Value = "'75.75"
$ws.Range("E22").Value = "  +2.87%  "
$ws.Range("D23").Value = "'0.549"
$ws.Range("E23").Value = "  -2.20%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -6.58%  "
$ws.Range("E26").Value = "  +6.40%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").Value = "'7.12"
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("D29").Value = "'7.83"
$ws.Range("E29").Value = "  -2.19%  "
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -6.17%  "
$ws.Range("D33").Value = "'22.92"
$ws.Range("E33").Value = "  -3.30%  "
$ws.Range("D34").Value = "'167.81"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D35").Value = "'6.81"
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("D36").Value = "'4.90"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("D37").Value = "3.399.33"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("D39").Value = "'0.0754"
$ws.Range("E39").Value = "  -2.86%  "
$ws.Range("D40").Value = "'25.30"
$ws.Range("E40").Value = "  -9.16%  "
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("E43").Value = "  -3.68%  "
$ws.Range("D44").Value = "2.458.79"
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("E45").Value = "  -2.56%  "
$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").Value = "'6.61"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").Value = "'22.10"
$ws.Range("E48").Value = "  -6.63%  "
$ws.Range("D49").Value = "'0.0257"
$ws.Range("E49").Value = "  -5.06%  "
$ws.Range("D50").Value = "'1.97"
$ws.Range("E50").Value = "  -4.72%  "
$ws.Range("E51").Value = "  -3.37%  "
